$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D cells to remain text-typed (values look numeric but source
# cells are stored as text in the original workbook).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.672.45"
$ws.Range("E2").Value = "  +0.55%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.591.00"
$ws.Range("E3").Value = "  -0.25%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "208.84"
$ws.Range("E5").Value = "  +0.97%  "
$ws.Range("E6").Value = "  +0.15%  "
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.27"
$ws.Range("E8").Value = "  +0.18%  "
$ws.Range("E9").Value = "  +0.18%  "
$ws.Range("E10").Value = "  +0.55%  "
$ws.Range("E11").Value = "  -0.16%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.818.67"
$ws.Range("E12").Value = "  -0.21%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.574.98"
$ws.Range("E13").Value = "  -1.17%  "
$ws.Range("E14").Value = "  -0.82%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.527"
$ws.Range("E15").Value = "  -2.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.702.49"
$ws.Range("E16").Value = "  +0.70%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.22"
$ws.Range("E17").Value = "  -0.15%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "217.70"
$ws.Range("E18").Value = "  +0.59%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0694"
$ws.Range("E19").Value = "  +0.65%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.34"
$ws.Range("E20").Value = "  -0.32%  "
$ws.Range("E21").Value = "  +0.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.15"
$ws.Range("E22").Value = "  -0.91%  "
$ws.Range("E23").Value = "  +0.84%  "
$ws.Range("E24").Value = "  -1.52%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.78"
$ws.Range("E25").Value = "  -0.57%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.95"
$ws.Range("E26").Value = "  +3.71%  "
$ws.Range("E27").Value = "  +0.18%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.08"
$ws.Range("E28").Value = "  +0.55%  "
$ws.Range("E29").Value = "  -0.16%  "
$ws.Range("E30").Value = "  -0.05%  "
$ws.Range("E31").Value = "  +1.71%  "
$ws.Range("E32").Value = "  -2.39%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.379.77"
$ws.Range("E33").Value = "  +1.81%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.97"
$ws.Range("E34").Value = "  +1.20%  "
$ws.Range("E35").Value = "  +0.13%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.964"
$ws.Range("E36").Value = "  +0.52%  "
$ws.Range("E37").Value = "  +0.21%  "
$ws.Range("E38").Value = "  +2.25%  "
$ws.Range("E39").Value = "  -0.45%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.826"
$ws.Range("E40").Value = "  +1.65%  "
$ws.Range("E41").Value = "  +0.22%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.982"
$ws.Range("E42").Value = "  +1.65%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "64.36"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.17"
$ws.Range("E44").Value = "  +4.32%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.75"
$ws.Range("E45").Value = "  +0.38%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.25"
$ws.Range("E46").Value = "  -1.12%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.729.19"
$ws.Range("E47").Value = "  -0.22%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.84"
$ws.Range("E48").Value = "  -1.29%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₇0999"
$ws.Range("E49").Value = "  +0.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0964"
$ws.Range("E50").Value = "  -0.43%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0495"
$ws.Range("E51").Value = "  -0.14%  "
